$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Power-flow results updated for the 380 kV case (rows 2-25, columns B:F and I:N)
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.047984526711065
$ws.Range("D2").Value = 1.04656188699557
$ws.Range("E2").Value = 1.061385451533013
$ws.Range("F2").Value = 1.068325505254562
$ws.Range("I2").Value = 1.040806349766163
$ws.Range("J2").Value = 1.053030688341007
$ws.Range("K2").Value = 1.049326976296874
$ws.Range("L2").Value = 1.06410965527328
$ws.Range("M2").Value = 1.071030989235665
$ws.Range("N2").Value = 1.021376121426462
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.049364455045103
$ws.Range("D3").Value = 1.047575659507991
$ws.Range("E3").Value = 1.062747371874247
$ws.Range("F3").Value = 1.069742620990987
$ws.Range("I3").Value = 1.041154290034212
$ws.Range("J3").Value = 1.054057197472237
$ws.Range("K3").Value = 1.050151891331896
$ws.Range("L3").Value = 1.065284807494826
$ws.Range("M3").Value = 1.072262567953305
$ws.Range("N3").Value = 1.021728940453259
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.050256353290962
$ws.Range("D4").Value = 1.048230502198343
$ws.Range("E4").Value = 1.063627962558272
$ws.Range("F4").Value = 1.070658939636541
$ws.Range("I4").Value = 1.041377288363388
$ws.Range("J4").Value = 1.054719963212003
$ws.Range("K4").Value = 1.050683925660687
$ws.Range("L4").Value = 1.06604400981479
$ws.Range("M4").Value = 1.073058303544114
$ws.Range("N4").Value = 1.021956457386525
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.050631071213899
$ws.Range("D5").Value = 1.048505528964791
$ws.Range("E5").Value = 1.063998008653675
$ws.Range("F5").Value = 1.071044009067763
$ws.Range("I5").Value = 1.041470525350647
$ws.Range("J5").Value = 1.054998245359728
$ws.Range("K5").Value = 1.050907178866133
$ws.Range("L5").Value = 1.066362895266312
$ws.Range("M5").Value = 1.073392552922058
$ws.Range("N5").Value = 1.02205191978018
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.050693974305846
$ws.Range("D6").Value = 1.04855169147101
$ws.Range("E6").Value = 1.064060132106528
$ws.Range("F6").Value = 1.071108655204873
$ws.Range("I6").Value = 1.041486150295108
$ws.Range("J6").Value = 1.055044950055241
$ws.Range("K6").Value = 1.050944639862228
$ws.Range("L6").Value = 1.066416421045967
$ws.Range("M6").Value = 1.073448658690854
$ws.Range("N6").Value = 1.022067937478276
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.050261361211527
$ws.Range("D7").Value = 1.048234178173336
$ws.Range("E7").Value = 1.063632907732794
$ws.Range("F7").Value = 1.070664085540364
$ws.Range("I7").Value = 1.041378536208796
$ws.Range("J7").Value = 1.054723682981156
$ws.Range("K7").Value = 1.05068691040441
$ws.Range("L7").Value = 1.066048271885624
$ws.Range("M7").Value = 1.073062770887284
$ws.Range("N7").Value = 1.021957733688055
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.048451091692931
$ws.Range("D8").Value = 1.046904733076458
$ws.Range("E8").Value = 1.061845858825428
$ws.Range("F8").Value = 1.068804563563898
$ws.Range("I8").Value = 1.040924382205356
$ws.Range("J8").Value = 1.053377905201996
$ws.Range("K8").Value = 1.049606122305749
$ws.Range("L8").Value = 1.064507054308982
$ws.Range("M8").Value = 1.07144745351495
$ws.Range("N8").Value = 1.02149552058148
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.045253217307409
$ws.Range("D9").Value = 1.044553253029979
$ws.Range("E9").Value = 1.058691573228852
$ws.Range("F9").Value = 1.065522659502527
$ws.Range("I9").Value = 1.040107635961389
$ws.Range("J9").Value = 1.05099517274108
$ws.Range("K9").Value = 1.047688179144339
$ws.Range("L9").Value = 1.061781862620274
$ws.Range("M9").Value = 1.06859184554627
$ws.Range("N9").Value = 1.020675012231528
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043115632428281
$ws.Range("D10").Value = 1.042979482235683
$ws.Range("E10").Value = 1.056584885391294
$ws.Range("F10").Value = 1.063330937953262
$ws.Range("I10").Value = 1.039551978053284
$ws.Range("J10").Value = 1.049398858440118
$ws.Range("K10").Value = 1.046400327892118
$ws.Range("L10").Value = 1.059958513605715
$ws.Range("M10").Value = 1.066681646567291
$ws.Range("N10").Value = 1.020123879519072
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.042188617174637
$ws.Range("D11").Value = 1.042296529774853
$ws.Range("E11").Value = 1.055671694042146
$ws.Range("F11").Value = 1.062380933646575
$ws.Range("I11").Value = 1.039308704895993
$ws.Range("J11").Value = 1.048705728982831
$ws.Range("K11").Value = 1.045840446931048
$ws.Range("L11").Value = 1.059167371222747
$ws.Range("M11").Value = 1.065852916986255
$ws.Range("N11").Value = 1.019884237518317
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.04184406169352
$ws.Range("D12").Value = 1.042042621899563
$ws.Range("E12").Value = 1.055332341063342
$ws.Range("F12").Value = 1.062027907504409
$ws.Range("I12").Value = 1.039217939527765
$ws.Range("J12").Value = 1.048447977627891
$ws.Range("K12").Value = 1.045632143231858
$ws.Range("L12").Value = 1.058873257249252
$ws.Range("M12").Value = 1.065544844147237
$ws.Range("N12").Value = 1.019795072419033
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.041917980105454
$ws.Range("D13").Value = 1.042097096419169
$ws.Range("E13").Value = 1.055405140445102
$ws.Range("F13").Value = 1.062103639793966
$ws.Range("I13").Value = 1.039237427267967
$ws.Range("J13").Value = 1.048503279464824
$ws.Range("K13").Value = 1.04567684047885
$ws.Range("L13").Value = 1.058936357019998
$ws.Range("M13").Value = 1.065610938023731
$ws.Range("N13").Value = 1.019814205509897
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.042160140651086
$ws.Range("D14").Value = 1.042275546361951
$ws.Range("E14").Value = 1.055643646173652
$ws.Range("F14").Value = 1.062351755535235
$ws.Range("I14").Value = 1.039301210431596
$ws.Range("J14").Value = 1.048684429162803
$ws.Range("K14").Value = 1.045823235417862
$ws.Range("L14").Value = 1.059143064767907
$ws.Range("M14").Value = 1.065827456633586
$ws.Range("N14").Value = 1.01987687020023
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.042309314314181
$ws.Range("D15").Value = 1.042385464816571
$ws.Range("E15").Value = 1.055790576976405
$ws.Range("F15").Value = 1.062504607482868
$ws.Range("I15").Value = 1.039340455903222
$ws.Range("J15").Value = 1.048796002613287
$ws.Range("K15").Value = 1.045913389147367
$ws.Range("L15").Value = 1.059270391157805
$ws.Range("M15").Value = 1.06596082816954
$ws.Range("N15").Value = 1.019915459879416
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043177124772141
$ws.Range("D16").Value = 1.043024775640004
$ws.Range("E16").Value = 1.056645469737099
$ws.Range("F16").Value = 1.06339396558221
$ws.Range("I16").Value = 1.039568066891987
$ws.Range("J16").Value = 1.049444818430102
$ws.Range("K16").Value = 1.046437438003926
$ws.Range("L16").Value = 1.060010984553864
$ws.Range("M16").Value = 1.066736612443955
$ws.Range("N16").Value = 1.020139762619844
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.043721092842832
$ws.Range("D17").Value = 1.043425394438793
$ws.Range("E17").Value = 1.057181454649142
$ws.Range("F17").Value = 1.063951571425792
$ws.Range("I17").Value = 1.039710125249486
$ws.Range("J17").Value = 1.049851287606436
$ws.Range("K17").Value = 1.046765559728417
$ws.Range("L17").Value = 1.060475101892121
$ws.Range("M17").Value = 1.067222809130278
$ws.Range("N17").Value = 1.020280193604293
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.044038242973413
$ws.Range("D18").Value = 1.043658924193603
$ws.Range("E18").Value = 1.057493991248176
$ws.Range("F18").Value = 1.064276719861759
$ws.Range("I18").Value = 1.03979272803006
$ws.Range("J18").Value = 1.050088189755138
$ws.Range("K18").Value = 1.04695673237917
$ws.Range("L18").Value = 1.060745657698297
$ws.Range("M18").Value = 1.067506245334535
$ws.Range("N18").Value = 1.020362008438758
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.044146359885079
$ws.Range("D19").Value = 1.043738527420258
$ws.Range("E19").Value = 1.057600542319161
$ws.Range("F19").Value = 1.06438757137246
$ws.Range("I19").Value = 1.039820849828798
$ws.Range("J19").Value = 1.050168936116964
$ws.Range("K19").Value = 1.047021880916816
$ws.Range("L19").Value = 1.060837883930436
$ws.Range("M19").Value = 1.067602863776481
$ws.Range("N19").Value = 1.020389888912837
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.043662744460169
$ws.Range("D20").Value = 1.043382426798064
$ws.Range("E20").Value = 1.057123958339666
$ws.Range("F20").Value = 1.063891755323237
$ws.Range("I20").Value = 1.039694910375687
$ws.Range("J20").Value = 1.049807696422659
$ws.Range("K20").Value = 1.046730377658326
$ws.Range("L20").Value = 1.060425322666255
$ws.Range("M20").Value = 1.067170660805717
$ws.Range("N20").Value = 1.020265136648396
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.042088836550934
$ws.Range("D21").Value = 1.042223003650544
$ws.Range("E21").Value = 1.055573416451141
$ws.Range("F21").Value = 1.062278695883398
$ws.Range("I21").Value = 1.03928243901587
$ws.Range("J21").Value = 1.048631093194841
$ws.Range("K21").Value = 1.045780135143064
$ws.Range("L21").Value = 1.059082201387803
$ws.Range("M21").Value = 1.065763704126613
$ws.Range("N21").Value = 1.019858421201747
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.041097977862293
$ws.Range("D22").Value = 1.041492701985143
$ws.Range("E22").Value = 1.054597640925948
$ws.Range("F22").Value = 1.06126361763142
$ws.Range("I22").Value = 1.039020770021493
$ws.Range("J22").Value = 1.047889623125874
$ws.Range("K22").Value = 1.045180716953751
$ws.Range("L22").Value = 1.058236287655506
$ws.Range("M22").Value = 1.064877670016102
$ws.Range("N22").Value = 1.019601826427349
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.041623374295065
$ws.Range("D23").Value = 1.041879975666694
$ws.Range("E23").Value = 1.055115004065105
$ws.Range("F23").Value = 1.061801815551216
$ws.Range("I23").Value = 1.039159707348517
$ws.Range("J23").Value = 1.048282852428198
$ws.Range("K23").Value = 1.045498667128729
$ws.Range("L23").Value = 1.058684860725774
$ws.Range("M23").Value = 1.065347510266153
$ws.Range("N23").Value = 1.019737935758888
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.043689110006344
$ws.Range("D24").Value = 1.043401842471415
$ws.Range("E24").Value = 1.057149938738149
$ws.Range("F24").Value = 1.063918783934133
$ws.Range("I24").Value = 1.039701786117634
$ws.Range("J24").Value = 1.049827393971867
$ws.Range("K24").Value = 1.046746275586132
$ws.Range("L24").Value = 1.060447816237356
$ws.Range("M24").Value = 1.067194224864694
$ws.Range("N24").Value = 1.020271940535758
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.046080919703083
$ws.Range("D25").Value = 1.045162231784301
$ws.Range("E25").Value = 1.059507686024743
$ws.Range("F25").Value = 1.066371755711521
$ws.Range("I25").Value = 1.040320745060844
$ws.Range("J25").Value = 1.051612528737148
$ws.Range("K25").Value = 1.048185626023314
$ws.Range("L25").Value = 1.062487526475378
$ws.Range("M25").Value = 1.06933120769131
$ws.Range("N25").Value = 1.020887855358813
